$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "mass" column header to "mass or m/z"
$ws.Range("B1").Value = "mass or m/z"

# Set column B width (approx. 12.46 characters wide)
$ws.Columns.Item(2).ColumnWidth = 11.6

# Update the view: zoom to 180% and move the selection to D12
$excel.ActiveWindow.Zoom = 180
$ws.Range("D12").Select()
